{"js": "// Target edit (see commit \"Add files via upload\"):\n//   1. \"Case-\u00c9\" -> \"Case- \u00c9\"   (a space is inserted between \"Case-\" and \"\u00c9\" so the\n//      word \"Case-\" is no longer glued to the following \"\u00c9 uma forma de reduzir...\")\n//   2. Remove bold formatting from the run\n//      \"O conte\u00fado de uma vari\u00e1vel \u00e9 comparado com um valor constante, e caso a\n//       compara\u00e7\u00e3o seja verdadeira, um determinado comando \u00e9 executado\"\n//   3. Remove bold formatting from the run\n//      \"a execu\u00e7\u00e3o das instru\u00e7\u00f5es vai continuar at\u00e9 que uma condi\u00e7\u00e3o seja verdadeira.\"\n\n// 1) Insert the missing space right after \"Case-\" (and before \"\u00c9\").\nconst caseDash = context.document.body.search(\"Case-\u00c9\", { matchCase: true });\ncaseDash.load(\"text\");\nawait context.sync();\n\nif (caseDash.items.length > 0) {\n  // \"Case-\u00c9\" was still glued together -> insert a space right after \"Case-\"\n  const afterCase = context.document.body.search(\"Case-\", { matchCase: true });\n  afterCase.load(\"text\");\n  await context.sync();\n  afterCase.items[0].insertText(\" \", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 2) & 3) Strip bold from the two specific runs of text.\nconst boldTargets = [\n  \"O conte\u00fado de uma vari\u00e1vel \u00e9 comparado com um valor constante, e caso a compara\u00e7\u00e3o seja verdadeira, um determinado comando \u00e9 executado\",\n  \"a execu\u00e7\u00e3o das instru\u00e7\u00f5es vai continuar at\u00e9 que uma condi\u00e7\u00e3o seja verdadeira.\"\n];\n\nfor (const phrase of boldTargets) {\n  const found = context.document.body.search(phrase, { matchCase: true });\n  found.load(\"text,font\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].font.bold = false;\n  }\n  await context.sync();\n}\n", "ps1": "# Target edit (see commit \"Add files via upload\"):\n#   1. \"Case-\u00c9\" -> \"Case- \u00c9\"   (insert a space right after \"Case-\" so the word is\n#      no longer glued to the following \"\u00c9 uma forma de reduzir...\")\n#   2. Remove bold formatting from the run:\n#      \"O conte\u00fado de uma vari\u00e1vel \u00e9 comparado com um valor constante, e caso a\n#       compara\u00e7\u00e3o seja verdadeira, um determinado comando \u00e9 executado\"\n#   3. Remove bold formatting from the run:\n#      \"a execu\u00e7\u00e3o das instru\u00e7\u00f5es vai continuar at\u00e9 que uma condi\u00e7\u00e3o seja verdadeira.\"\n\n$d = $word.ActiveDocument\n\n# 1) Insert the missing space right after \"Case-\" (and before \"\u00c9\").\n$caseRange = $d.Content\n$caseRange.Find.MatchCase = $true\n$caseRange.Find.Text = \"Case-\"\nif ($caseRange.Find.Execute()) {\n    $caseRange.InsertAfter(\" \")\n}\n\n# 2) Strip bold from \"O conte\u00fado de uma vari\u00e1vel ... \u00e9 executado\"\n$r2 = $d.Content\n$r2.Find.MatchCase = $true\n$r2.Find.Text = \"O conte\u00fado de uma vari\u00e1vel \u00e9 comparado com um valor constante, e caso a compara\u00e7\u00e3o seja verdadeira, um determinado comando \u00e9 executado\"\nif ($r2.Find.Execute()) {\n    $r2.Font.Bold = 0\n    $r2.Font.BoldBi = 0\n}\n\n# 3) Strip bold from \"a execu\u00e7\u00e3o das instru\u00e7\u00f5es vai continuar at\u00e9 que uma condi\u00e7\u00e3o seja verdadeira.\"\n$r3 = $d.Content\n$r3.Find.MatchCase = $true\n$r3.Find.Text = \"a execu\u00e7\u00e3o das instru\u00e7\u00f5es vai continuar at\u00e9 que uma condi\u00e7\u00e3o seja verdadeira.\"\nif ($r3.Find.Execute()) {\n    $r3.Font.Bold = 0\n    $r3.Font.BoldBi = 0\n}\n"}
